$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, matching the style of the other headers (G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H10 with 0 (no save), H11 with 1 (saved)
$ws.Range("H2:H10").Value = 0
$ws.Range("H11").Value = 1
